{"js": "// Load the existing paragraphs in the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Insert a new, empty paragraph in front of the first paragraph\n//    (\"Unify correlation and variance nomenclature\"). The new paragraph\n//    keeps the same visual (hanging) indent the bulleted list uses\n//    (720 left / 360 hanging twips = 36pt / -18pt) but is not part of\n//    the numbered/bulleted list itself, so it carries no list style or\n//    numbering.\nconst firstParagraph = paragraphs.items[0];\nconst blankParagraph = firstParagraph.insertParagraph(\"\", Word.InsertLocation.before);\nblankParagraph.style = \"Normal\";\nblankParagraph.leftIndent = 36;       // -> w:ind w:left=\"720\"\nblankParagraph.firstLineIndent = -18; // -> w:ind w:hanging=\"360\"\n\n// 2) The last paragraph in the list (currently empty) gets its text filled in.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertText(\n  \"Base variables selection in significance of effects\",\n  Word.InsertLocation.end\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Insert a new, empty paragraph in front of the first paragraph\n#    (\"Unify correlation and variance nomenclature\"). The new paragraph\n#    keeps the same visual (hanging) indent the bulleted list uses\n#    (720 left / 360 hanging twips = 36pt / -18pt) but is not itself part\n#    of the numbered/bulleted list, so it carries no list style/numbering.\n$firstPara = $d.Paragraphs(1)\n$insertionRange = $firstPara.Range\n$insertionRange.Collapse(1)            # wdCollapseStart\n$insertionRange.InsertParagraphBefore()\n\n$blankPara = $d.Paragraphs(1)\n$blankPara.Range.Style = \"Normal\"\n$blankPara.Range.ParagraphFormat.LeftIndent = 36        # -> w:ind w:left=\"720\"\n$blankPara.Range.ParagraphFormat.FirstLineIndent = -18  # -> w:ind w:hanging=\"360\"\n\n# 2) The last paragraph in the list (currently empty) gets its text filled in.\n$lastPara = $d.Paragraphs($d.Paragraphs.Count)\n$lastPara.Range.InsertAfter(\"Base variables selection in significance of effects\")\n"}
